$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the table on the sheet and add a new row to it so the table
# range (and autofilter) grows from A1:F74 to A1:F75.
$tbl = $ws.ListObjects.Item(1)
[void]$tbl.ListRows.Add()

# Copy formatting (number format / style) from the last existing data
# row (74) down into the newly created row (75) so the new cells keep
# the same look (date format in column A, centered numbers elsewhere).
$ws.Range("A74:F74").Copy()
$ws.Range("A75:F75").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row with the latest data point.
$ws.Range("A75").Value = 43977
$ws.Range("B75").Value = 552
$ws.Range("C75").Value = 212
$ws.Range("D75").Value = 405
$ws.Range("E75").Value = 28
$ws.Range("F75").Value = 23

# Update the selection to match where the user ended up after entering
# the new data.
$ws.Activate()
[void]$ws.Range("C75").Select()
